# readxl test fixture "blanks.xlsx" rework:
#  - rename the lone sheet to "different_rows" (data/layout unchanged)
#  - add two more sheets ("same_row_first", "same_row_middle") that reuse the
#    same header row, but skip a row somewhere in the middle (blank row),
#    exercising "skipping and/or blank rows" cell ingest.

$wb = $excel.ActiveWorkbook

# --- sheet 1: just rename, keep its existing data (A1=x, B1=y, B2=a, A3=1) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "different_rows"

# --- sheet 2: header row, then a completely blank row 2, then data row 3 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "same_row_first"
$ws2.Range("A1").Value = "x"
$ws2.Range("B1").Value = "y"
$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "a"

# --- sheet 3: header row, data row 2, blank row 3, data row 4 ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "same_row_middle"
$ws3.Range("A1").Value = "x"
$ws3.Range("B1").Value = "y"
$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = "a"
$ws3.Range("A4").Value = 2
$ws3.Range("B4").Value = "b"

# --- selections matching each sheet's saved view ---
[void]$ws1.Range("A1:B3").Select()
[void]$ws2.Range("A1:B3").Select()
[void]$ws3.Range("A5").Select()

# "same_row_first" is the tab that ends up active/selected in the workbook
$ws2.Activate()
